# Refresh the "想去人数" (want-to-go count, column F) and, in one case,
# the "最低票价" (min price, column G) stats across the four sheets
# (展览 / 演出 / 本地生活 / 全部类型) to match the latest scrape.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 2306
$ws1.Cells.Item(3, 6).Value = 383
$ws1.Cells.Item(4, 6).Value = 193
$ws1.Cells.Item(5, 6).Value = 306
$ws1.Cells.Item(6, 6).Value = 486
$ws1.Cells.Item(8, 6).Value = 730
$ws1.Cells.Item(9, 6).Value = 527
$ws1.Cells.Item(10, 6).Value = 712
$ws1.Cells.Item(11, 6).Value = 376
$ws1.Cells.Item(12, 6).Value = 73
$ws1.Cells.Item(13, 6).Value = 374
$ws1.Cells.Item(15, 6).Value = 989
$ws1.Cells.Item(16, 6).Value = 17210
$ws1.Cells.Item(17, 6).Value = 371
$ws1.Cells.Item(18, 6).Value = 39
$ws1.Cells.Item(19, 6).Value = 172
$ws1.Cells.Item(20, 6).Value = 269
$ws1.Cells.Item(21, 6).Value = 159
$ws1.Cells.Item(22, 6).Value = 127
$ws1.Cells.Item(23, 6).Value = 11
$ws1.Cells.Item(24, 6).Value = 151
$ws1.Cells.Item(26, 6).Value = 297
$ws1.Cells.Item(27, 6).Value = 122

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(4, 6).Value = 164
$ws2.Cells.Item(6, 6).Value = 190
$ws2.Cells.Item(8, 6).Value = 3313
$ws2.Cells.Item(10, 6).Value = 34
$ws2.Cells.Item(16, 6).Value = 2694

$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(3, 6).Value = 68
$ws3.Cells.Item(4, 6).Value = 504
$ws3.Cells.Item(5, 6).Value = 190

$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 68
$ws4.Cells.Item(6, 6).Value = 2306
$ws4.Cells.Item(7, 6).Value = 504
$ws4.Cells.Item(8, 6).Value = 383
$ws4.Cells.Item(9, 6).Value = 193
$ws4.Cells.Item(10, 6).Value = 306
$ws4.Cells.Item(11, 6).Value = 486
$ws4.Cells.Item(12, 6).Value = 164
$ws4.Cells.Item(15, 6).Value = 190
$ws4.Cells.Item(16, 6).Value = 190
$ws4.Cells.Item(17, 6).Value = 730
$ws4.Cells.Item(18, 6).Value = 527
$ws4.Cells.Item(19, 6).Value = 712
$ws4.Cells.Item(20, 6).Value = 376
$ws4.Cells.Item(21, 6).Value = 73
$ws4.Cells.Item(22, 6).Value = 374
$ws4.Cells.Item(24, 6).Value = 989
$ws4.Cells.Item(27, 6).Value = 3313
$ws4.Cells.Item(29, 6).Value = 34
$ws4.Cells.Item(31, 6).Value = 371
$ws4.Cells.Item(32, 6).Value = 39
$ws4.Cells.Item(33, 6).Value = 172
$ws4.Cells.Item(36, 6).Value = 269
$ws4.Cells.Item(37, 6).Value = 159
$ws4.Cells.Item(38, 6).Value = 127
$ws4.Cells.Item(39, 6).Value = 11
$ws4.Cells.Item(42, 6).Value = 151
$ws4.Cells.Item(44, 6).Value = 297
$ws4.Cells.Item(45, 6).Value = 122
$ws4.Cells.Item(46, 6).Value = 2694

# Special case: sheet4 (全部类型) row 25 - item now has a price instead of "temporarily sold out"
$ws4.Cells.Item(25, 6).Value = 17211
$ws4.Cells.Item(25, 7).Value = 85
